$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "PREPARATION" cell (F2) text: No. Urut changes from 1369 to 2962
$ws.Range("F2").Value = "Username : 31246;`nPassword : bni1234;`nRole : 20/21 - Analis Investasi/Asisten Investasi;`nNo. Urut : 2962"

# Update the "NO_URUT" cell (M2) value from 1369 to 2962
$ws.Range("M2").Value = 2962

# Update the active selection to G2 (from M1)
$ws.Range("G2").Select()
